$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.068.19"
$ws.Range("E2").Value = "  -2.54%  "

$ws.Range("D3").Value = "2.360.36"
$ws.Range("E3").Value = "  -3.63%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'499.44"
$ws.Range("E5").Value = "  -2.03%  "

$ws.Range("D6").Value = "'128.46"
$ws.Range("E6").Value = "  -3.49%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D8").Value = "'0.544"
$ws.Range("E8").Value = "  -2.29%  "

$ws.Range("D9").Value = "2.362.11"
$ws.Range("E9").Value = "  -3.51%  "

$ws.Range("D10").Value = "'0.0974"
$ws.Range("E10").Value = "  -0.50%  "

$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("D12").Value = "'4.73"
$ws.Range("E12").Value = "  +3.06%  "

$ws.Range("D13").Value = "'0.321"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").Value = "2.780.83"
$ws.Range("E14").Value = "  -3.50%  "

$ws.Range("D15").Value = "56.050.78"
$ws.Range("E15").Value = "  -2.56%  "

$ws.Range("D16").Value = "'21.32"
$ws.Range("E16").Value = "  -2.34%  "

$ws.Range("D17").Value = "'0.0000131"
$ws.Range("E17").Value = "  -1.56%  "

$ws.Range("D18").Value = "2.422.12"
$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("D19").Value = "'9.96"
$ws.Range("E19").Value = "  -3.22%  "

$ws.Range("D20").Value = "'4.01"
$ws.Range("E20").Value = "  -2.29%  "

$ws.Range("D21").Value = "'305.56"
$ws.Range("E21").Value = "  -2.65%  "

$ws.Range("D22").Value = "'6.24"
$ws.Range("E22").Value = "  -2.94%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'65.16"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "'1.01"
$ws.Range("E25").Value = "  +1.33%  "

$ws.Range("D26").Value = "'0.368"
$ws.Range("E26").Value = "  -3.38%  "

$ws.Range("E27").Value = "  -5.57%  "

$ws.Range("D28").Value = "'7.17"
$ws.Range("E28").Value = "  -4.82%  "

$ws.Range("D29").Value = "'171.11"
$ws.Range("E29").Value = "  -1.38%  "

$ws.Range("D30").Value = "0.0₃0706"
$ws.Range("E30").Value = "  -3.58%  "

$ws.Range("D31").Value = "'1.63"
$ws.Range("E31").Value = "  -3.52%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.42%  "

$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.07"
$ws.Range("E34").Value = "  -5.08%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'5.70"
$ws.Range("E35").Value = "  -7.45%  "

$ws.Range("D36").Value = "'17.51"
$ws.Range("E36").Value = "  -2.64%  "

$ws.Range("E37").Value = "  -5.93%  "

$ws.Range("D38").Value = "'3.72"
$ws.Range("E38").Value = "  -2.75%  "

$ws.Range("D39").Value = "'36.01"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").Value = "'0.784"
$ws.Range("E40").Value = "  -3.44%  "

$ws.Range("D41").Value = "'1.37"
$ws.Range("E41").Value = "  -6.11%  "

$ws.Range("D42").Value = "'129.00"
$ws.Range("E42").Value = "  -5.08%  "

$ws.Range("D43").Value = "'3.33"
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("D44").Value = "'4.68"
$ws.Range("E44").Value = "  -6.35%  "

$ws.Range("D45").Value = "'0.559"
$ws.Range("E45").Value = "  -2.27%  "

$ws.Range("D46").Value = "'0.0901"
$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("D47").Value = "'238.48"
$ws.Range("E47").Value = "  -6.76%  "

$ws.Range("D48").Value = "'0.0478"
$ws.Range("E48").Value = "  -2.75%  "

$ws.Range("D49").Value = "'0.0206"
$ws.Range("E49").Value = "  -3.48%  "

$ws.Range("D50").Value = "'16.90"
$ws.Range("E50").Value = "  -0.81%  "

$ws.Range("E51").Value = "  -0.68%  "
